$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("A17").Value = "Complete le chargement des niveaux"
$ws.Range("B17").Value = "30min"
$ws.Range("C17").Value = "TP4"
$ws.Range("F17").Value = "OK"

$ws.Range("A18").Value = "Gére les collisions entre les caisses et les murs"
$ws.Range("B18").Value = "30min"
$ws.Range("C18").Value = "TP4"
$ws.Range("F18").Value = "OK"

$ws.Range("A19").Value = "Améliore la map"
$ws.Range("B19").Value = "30min"
$ws.Range("C19").Value = "TP4"
$ws.Range("F19").Value = "OK"

$ws.Range("F22").Select()
